$wb = $excel.ActiveWorkbook

# Rename the "payroll_dec" sheet to "payroll_feb"
$sheet = $wb.Worksheets.Item("payroll_dec")
$sheet.Name = "payroll_feb"

# Make the renamed sheet the active/selected sheet
$sheet.Activate()
